# Update the "dSF" (column F) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 10
    9  = -4
    10 = -6
    15 = 6
    16 = 1
    17 = -3
    20 = -5
    21 = 7
    24 = -6
    26 = 5
    28 = -4
    30 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
